# Applies the diff:
#  - removes the custom paragraph style "AbstractTitle" ("Abstract Title")
#  - changes the "Abstract" style's paragraph spacing-before from 100 (5pt) to 300 (15pt)

$d = $word.ActiveDocument

# Remove the "AbstractTitle" style entirely.
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Delete()

# Update the "Abstract" style's space-before (twips 100 -> 300, i.e. points 5 -> 15).
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15
